$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, taken from the commit diff.
# Columns B/C (Coin name, Link) are plain text.
# Columns D/E (Price, Volume) look numeric/percent, so Excel would
# auto-convert them; they are written with a leading apostrophe so
# they land as literal text (matching the original inlineStr cells),
# then the cell style is reset to 'Normal' so the quote-prefix flag
# does not linger as a visible style change.
$updates = [ordered]@{
    'D2' = '293.22'
    'E2' = '0.03%'
    'D3' = '40.46'
    'E3' = '1.60%'
    'D4' = '5.007'
    'E4' = '-0.59%'
    'D5' = '0.07347'
    'E5' = '-0.45%'
    'B6' = 'FTXToken'
    'C6' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D6' = '1.571'
    'E6' = '0.13%'
    'B7' = 'MXToken'
    'C7' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D7' = '0.9253'
    'E7' = '0.51%'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D8' = '2.362'
    'E8' = '-1.55%'
    'B9' = 'LiechtensteinCryptoassetsExchange'
    'C9' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D9' = '0.1164'
    'E9' = '0.29%'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1780'
    'E10' = '1.85%'
    'B11' = 'BitrueCoin'
    'C11' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D11' = '0.04373'
    'E11' = '5.24%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.08717'
    'E12' = '0.65%'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.1053'
    'E13' = '0.15%'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001271'
    'E14' = '1.46%'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D15' = '0.005888'
    'E15' = '-1.06%'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D16' = '3.351'
    'E16' = '-0.13%'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D17' = '4.296'
    'E17' = '-0.37%'
    'B18' = 'BitpandaEcosystemToken'
    'C18' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D18' = '0.3307'
    'E18' = '-0.26%'
    'B19' = 'MCDex'
    'C19' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D19' = '7.803'
    'E19' = '2.93%'
    'B20' = 'ProBitToken'
    'C20' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D20' = '0.1390'
    'E20' = '2.18%'
    'B21' = 'ZBToken'
    'C21' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
    'D21' = '0.2769'
    'E21' = '-1.89%'
    'B22' = 'CoinExToken'
    'C22' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D22' = '0.03924'
    'E22' = '2.72%'
    'B23' = 'BitKan'
    'C23' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'D23' = '0.001262'
    'E23' = '-2.24%'
    'B24' = 'HotbitToken'
    'C24' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'D24' = '0.003670'
    'E24' = '1.79%'
    'D25' = '0.0001200'
    'E25' = '-8.29%'
    'D26' = '0.0003721'
    'E26' = '-0.78%'
    'D38' = '0.02331'
    'E38' = '1.38%'
    'D39' = '0.05091'
    'E39' = '2.11%'
    'D40' = '0.005542'
    'E40' = '38.00%'
    'D41' = '0.007860'
    'E41' = '1.54%'
    'D42' = '0.1293'
    'E42' = '1.67%'
    'D43' = '0.007377'
    'E43' = '-0.56%'
    'D44' = '0.008050'
    'E44' = '1.29%'
    'D45' = '0.2915'
    'E45' = '-8.04%'
    'D46' = '0.00006235'
    'E46' = '-3.69%'
    'E47' = '-0.79%'
    'D48' = '0.04749'
    'E48' = '-81.14%'
    'E49' = '-0.79%'
    'E50' = '-0.79%'
}

$textLikeColumns = @('D', 'E')

foreach ($addr in $updates.Keys) {
    $col = [regex]::Match($addr, '^[A-Z]+').Value
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($textLikeColumns -contains $col) {
        $cell.Value = "'" + $value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $value
    }
}
